# Update the SQL Saturday Colorado Springs 2022 event checklist:
# - new event date / time
# - add registration URL and contact-email hyperlinks
# - adjust the "Value" column width so the new URL fits
# - move the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("New Event Checklist")

# Event Date: Mar 21, 2022 -> Mar 19, 2022
$ws.Range("C7").Value = "Mar 19, 2022"

# Event Time: 1000 - 1630 MST -> 1000 - 1730 MST
$ws.Range("C8").Value = "1000 - 1730 MST"

# Registration URL (was blank) -> Eventbrite link, shown as a hyperlink
$ws.Hyperlinks.Add(
    $ws.Range("C11"),
    "https://www.eventbrite.com/e/sql-saturday-colorado-springs-tickets-247770376867",
    [Type]::Missing,
    [Type]::Missing,
    "https://www.eventbrite.com/e/sql-saturday-colorado-springs-tickets-247770376867"
) | Out-Null

# Contact Email (was blank) -> springssql@gmail.com, shown as a mailto hyperlink
$ws.Hyperlinks.Add(
    $ws.Range("C15"),
    "mailto:springssql@gmail.com",
    [Type]::Missing,
    [Type]::Missing,
    "springssql@gmail.com"
) | Out-Null

# Widen column C to fit the long registration URL (best-fit/auto-fit)
$ws.Columns.Item(3).EntireColumn.AutoFit() | Out-Null
$ws.Columns.Item(3).ColumnWidth = 145.92

# Move the active selection to C13
$ws.Range("C13").Select() | Out-Null
